# Generate Report for Handoff
# Updates the localization status report: marks items as "Ready for handoff",
# refreshes the handoff timestamps, and tightens the datetime column widths.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

# --- Status text: "Handed back: in sync with en-US" -> "Ready for handoff"
$overview.Range("E2").Value = "Ready for handoff"
$overview.Range("F2").Value = "Ready for handoff"
$zhcn.Range("C2").Value = "Ready for handoff"
$dede.Range("C2").Value = "Ready for handoff"

# --- Refresh "Latest HO Xliff Generate Date" / "Latest Handoff Datetime" timestamps
$overview.Range("G2").Value = "2016-09-06 23:13:49"
$dede.Range("H2").Value = "2016-09-06 23:13:49"
$zhcn.Range("H2").Value = "2016-09-06 23:13:44"

# --- Tighten datetime column widths (29.9777047293527 -> 17.2159881591797)
$overview.Columns.Item(5).ColumnWidth = 16.3
$overview.Columns.Item(6).ColumnWidth = 16.3
$zhcn.Columns.Item(3).ColumnWidth = 16.3
$dede.Columns.Item(3).ColumnWidth = 16.3
